$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the PLP_Page sheet (this is how PDP_Page was created - it keeps the
# same column widths / styles / page setup as the original sheet) and place it
# right after PLP_Page.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "PDP_Page"

# Wipe out the old PLP content in rows 2-9 (the product grid / filter rows plus
# the "collor filters" sub-heading) so we can fill in the new PDP locators.
$newSheet.Range("A2:D9").Clear()

# The "price range filters" heading row (old row 16) is removed entirely, which
# shifts the remaining price-filter rows up by one.
$newSheet.Rows.Item(16).Delete()

# New PDP_Page locator rows.
$newSheet.Range("A2").Value = "product_info"
$newSheet.Range("B2").Value = "class name"
$newSheet.Range("C2").Value = "AddProductToCartFormProductSummary_productSummary_row__3pIGN"

$newSheet.Range("A3").Value = "quantity"
$newSheet.Range("B3").Value = "class name"

$newSheet.Range("A4").Value = "quantity"
$newSheet.Range("B4").Value = "id"
$newSheet.Range("C4").Value = "otCSfuRyXDDgPcWMS82WHWpQ"

$newSheet.Range("C3").Value = "pt-3.5"

$newSheet.Range("A5").Value = "quantity"
$newSheet.Range("B5").Value = "xpath"
$newSheet.Range("C5").Value = '//*[@id="Flu6humgibSbVSKMGLipUKXy"]'

# The leftover styled-but-empty rows from the copied sheet (colour filter rows
# and price filter rows) keep their formatting but lose their values.
$newSheet.Range("A10:D15").ClearContents()
$newSheet.Range("A16:D21").ClearContents()

$newSheet.Range("C6").Select()
